# feat: add 2022-Q1 data
#
# The workbook gains a new per-quarter sheet "2022-Q1" (inserted between
# "2021-Q2" and "总计", reusing the same layout/content shape as
# "2021-Q2") and the "总计" summary sheet gets a new leading row with the
# 2022-Q1 totals (newest quarter first).

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text, $formatDonor) {
    # Force a numeric-looking string (e.g. "000593", "1.93") to be stored
    # as literal text rather than being auto-parsed into a number. The
    # leading apostrophe also stamps a "quote prefix" format onto the
    # cell, so immediately re-paste the (plain) formatting from a donor
    # cell to restore the formatting the cell is actually supposed to
    # have.
    $range.Value = "'" + $text
    $formatDonor.Copy()
    $range.PasteSpecial(-4122)
}

# ---------------------------------------------------------------------
# 0) duplicate the existing "总计" sheet (before touching it) so the
#    brand-new "总计" sheet keeps the same sheet-level formatting
#    (margins, outline properties, …) as the rest of the workbook
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")
$oldTotal.Copy($null, $oldTotal)
$total = $wb.Worksheets.Item("总计 (2)")

# ---------------------------------------------------------------------
# 1) old "总计" sheet -> becomes the new "2022-Q1" sheet
# ---------------------------------------------------------------------
$q1 = $oldTotal
$q1.Name = "2022-Q1"

# extend the header-row / index-column formatting (copied from the cells
# that already carry it) before writing the new values
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$q1.Range("A3:A5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $q1.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "000593", "易方达标普全球高端消费品指数增强(QDII)-美元现汇", "1.93", "92.46", "7.87", "0.1519", 5),
    @(1, "005676", "易方达标普全球高端消费品指数增强C(QDII) - 人民币", "1.93", "92.46", "7.87", "0.1519", 5),
    @(2, "118002", "易方达标普全球高端消费品指数增强A(QDII) - 人民币", "1.93", "92.46", "7.87", "0.1519", 5),
    @(3, "513080", "华安法国CAC40ETF（QDII）", "0.60", "96.69", "3.17", "0.0190", 10)
)

foreach ($row in $rows) {
    $r = [int]$row[0] + 2
    $q1.Cells.Item($r, 1).Value = $row[0]
    $q1.Cells.Item($r, 3).Value = $row[2]
    $donor = $q1.Cells.Item($r, 3)
    Set-TextValue $q1.Cells.Item($r, 2) $row[1] $donor
    Set-TextValue $q1.Cells.Item($r, 4) $row[3] $donor
    Set-TextValue $q1.Cells.Item($r, 5) $row[4] $donor
    Set-TextValue $q1.Cells.Item($r, 6) $row[5] $donor
    Set-TextValue $q1.Cells.Item($r, 7) $row[6] $donor
    $q1.Cells.Item($r, 8).Value = $row[7]
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) the duplicated sheet becomes the new "总计" sheet, rewritten with
#    the 2022-Q1 quarter on top and the pre-existing 2021-Q2 row below
# ---------------------------------------------------------------------
$total.Name = "总计"

# A2 already carries the index-column style (it was copied along with
# the rest of the old "总计" sheet); extend it down to the new A3 row
$total.Range("A2").Copy()
$total.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$total.Range("A3").Value = 1

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 4
$total.Cells.Item(2, 4).Value = 0.47

$total.Cells.Item(3, 2).Value = "2021-Q2"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.48
